$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-18 from serial date 45186 (2023-09-17)
# to serial date 45188 (2023-09-19).
$ws.Range("C2:C18").Value = 45188
